$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '23.491.97'
$ws.Range("E2").Value = '  +0.06%  '

Set-TextValue $ws.Range("D3") '1.639.95'
$ws.Range("E3").Value = '  -0.11%  '

Set-TextValue $ws.Range("D4") '0.9998'
$ws.Range("E4").Value = '  -0.40%  '

Set-TextValue $ws.Range("D5") '1.000'
$ws.Range("E5").Value = '  -0.54%  '

Set-TextValue $ws.Range("D6") '304.28'
$ws.Range("E6").Value = '  +0.06%  '

Set-TextValue $ws.Range("D7") '0.3784'
$ws.Range("E7").Value = '  +0.53%  '

Set-TextValue $ws.Range("D8") '51.72'
$ws.Range("E8").Value = '  -1.33%  '

Set-TextValue $ws.Range("D9") '0.3624'
$ws.Range("E9").Value = '  -0.61%  '

Set-TextValue $ws.Range("D10") '0.08172'
$ws.Range("E10").Value = '  +0.36%  '

Set-TextValue $ws.Range("D11") '1.232'
$ws.Range("E11").Value = '  -1.59%  '

Set-TextValue $ws.Range("D12") '1.001'
$ws.Range("E12").Value = '  -0.26%  '

Set-TextValue $ws.Range("D13") '22.52'
$ws.Range("E13").Value = '  -1.93%  '

Set-TextValue $ws.Range("D14") '6.465'
$ws.Range("E14").Value = '  -2.89%  '

Set-TextValue $ws.Range("D15") '7.370'
$ws.Range("E15").Value = '  +0.85%  '

Set-TextValue $ws.Range("D16") '0.00001239'
$ws.Range("E16").Value = '  -1.81%  '

Set-TextValue $ws.Range("D17") '1.634.79'
$ws.Range("E17").Value = '  -0.64%  '

Set-TextValue $ws.Range("D18") '95.55'
$ws.Range("E18").Value = '  +1.24%  '

Set-TextValue $ws.Range("D19") '0.06942'
$ws.Range("E19").Value = '  +0.00%  '

Set-TextValue $ws.Range("D20") '6.588'
$ws.Range("E20").Value = '  +0.33%  '

Set-TextValue $ws.Range("D21") '17.52'
$ws.Range("E21").Value = '  -3.72%  '

Set-TextValue $ws.Range("D22") '0.9996'
$ws.Range("E22").Value = '  -0.71%  '

Set-TextValue $ws.Range("D23") '12.51'
$ws.Range("E23").Value = '  -3.01%  '

Set-TextValue $ws.Range("D24") '23.496.24'
$ws.Range("E24").Value = '  +0.03%  '

Set-TextValue $ws.Range("D25") '2.509'
$ws.Range("E25").Value = '  +3.00%  '

Set-TextValue $ws.Range("D26") '3.049'
$ws.Range("E26").Value = '  -5.85%  '

Set-TextValue $ws.Range("D27") '21.16'
$ws.Range("E27").Value = '  -0.57%  '

Set-TextValue $ws.Range("D28") '152.11'
$ws.Range("E28").Value = '  +0.59%  '

Set-TextValue $ws.Range("D29") '5.250'
$ws.Range("E29").Value = '  -1.05%  '

Set-TextValue $ws.Range("D30") '132.83'
$ws.Range("E30").Value = '  -2.51%  '

Set-TextValue $ws.Range("D31") '1.817.21'
$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("B32").Value = 'WEMIXTOKEN'
$ws.Range("C32").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D32") '2.164'
$ws.Range("E32").Value = '  -6.99%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D33") '6.600'
$ws.Range("E33").Value = '  -4.53%  '

Set-TextValue $ws.Range("D34") '1.078'
$ws.Range("E34").Value = '  +11.88%  '

Set-TextValue $ws.Range("D35") '11.54'
$ws.Range("E35").Value = '  +5.51%  '

Set-TextValue $ws.Range("D36") '0.02761'
$ws.Range("E36").Value = '  -3.99%  '

$ws.Range("E37").Value = '  -2.86%  '

Set-TextValue $ws.Range("D38") '0.08780'
$ws.Range("E38").Value = '  -0.89%  '

Set-TextValue $ws.Range("D39") '0.07098'
$ws.Range("E39").Value = '  -2.51%  '

Set-TextValue $ws.Range("D40") '5.987'
$ws.Range("E40").Value = '  -4.86%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D41") '1.352'
$ws.Range("E41").Value = '  -1.86%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D42") '0.7048'
$ws.Range("E42").Value = '  -1.57%  '

Set-TextValue $ws.Range("D43") '12.16'
$ws.Range("E43").Value = '  -3.38%  '

Set-TextValue $ws.Range("D44") '15.65'
$ws.Range("E44").Value = '  -5.06%  '

Set-TextValue $ws.Range("D45") '0.6526'
$ws.Range("E45").Value = '  -0.98%  '

Set-TextValue $ws.Range("D46") '0.9995'
$ws.Range("E46").Value = '  -0.54%  '

Set-TextValue $ws.Range("D47") '2.276'
$ws.Range("E47").Value = '  -3.98%  '

$ws.Range("E48").Value = '  -0.78%  '

Set-TextValue $ws.Range("D49") '0.07983'
$ws.Range("E49").Value = '  -0.41%  '

Set-TextValue $ws.Range("D50") '129.00'
$ws.Range("E50").Value = '  +0.68%  '

Set-TextValue $ws.Range("D51") '1.193'
$ws.Range("E51").Value = '  -1.83%  '
